$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: set value, then copy E1 formatting (bold/border/center style)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Data cells F2:F67 -> "time_taken" timestamps recorded for each panel row
$timeTaken = @(
    "2021-10-05 13:39:29.640218",
    "2021-10-05 13:39:29.640229",
    "2021-10-05 13:39:29.640232",
    "2021-10-05 13:39:29.640235",
    "2021-10-05 13:39:29.640237",
    "2021-10-05 13:39:29.640240",
    "2021-10-05 13:39:29.640243",
    "2021-10-05 13:39:29.640245",
    "2021-10-05 13:39:29.640248",
    "2021-10-05 13:39:29.640250",
    "2021-10-05 13:39:29.640253",
    "2021-10-05 13:39:29.640255",
    "2021-10-05 13:39:29.640258",
    "2021-10-05 13:39:29.640260",
    "2021-10-05 13:39:29.640262",
    "2021-10-05 13:39:29.640265",
    "2021-10-05 13:39:29.640268",
    "2021-10-05 13:39:29.640270",
    "2021-10-05 13:39:29.640273",
    "2021-10-05 13:39:29.640275",
    "2021-10-05 13:39:29.640278",
    "2021-10-05 13:39:29.640280",
    "2021-10-05 13:39:29.640283",
    "2021-10-05 13:39:29.640285",
    "2021-10-05 13:39:29.640288",
    "2021-10-05 13:39:29.640291",
    "2021-10-05 13:39:29.640293",
    "2021-10-05 13:39:29.640296",
    "2021-10-05 13:39:29.640298",
    "2021-10-05 13:39:29.640301",
    "2021-10-05 13:39:29.640303",
    "2021-10-05 13:39:29.640306",
    "2021-10-05 13:39:29.640308",
    "2021-10-05 13:39:29.640311",
    "2021-10-05 13:39:29.640313",
    "2021-10-05 13:39:29.640316",
    "2021-10-05 13:39:29.640318",
    "2021-10-05 13:39:29.640321",
    "2021-10-05 13:39:29.640323",
    "2021-10-05 13:39:29.640326",
    "2021-10-05 13:39:29.640329",
    "2021-10-05 13:39:29.640331",
    "2021-10-05 13:39:29.640334",
    "2021-10-05 13:39:29.640336",
    "2021-10-05 13:39:29.640338",
    "2021-10-05 13:39:29.640341",
    "2021-10-05 13:39:29.640343",
    "2021-10-05 13:39:29.640346",
    "2021-10-05 13:39:29.640348",
    "2021-10-05 13:39:29.640351",
    "2021-10-05 13:39:29.640353",
    "2021-10-05 13:39:29.640355",
    "2021-10-05 13:39:29.640358",
    "2021-10-05 13:39:29.640361",
    "2021-10-05 13:39:29.640363",
    "2021-10-05 13:39:29.640366",
    "2021-10-05 13:39:29.640368",
    "2021-10-05 13:39:29.640371",
    "2021-10-05 13:39:29.640373",
    "2021-10-05 13:39:29.640375",
    "2021-10-05 13:39:29.640378",
    "2021-10-05 13:39:29.640380",
    "2021-10-05 13:39:29.640383",
    "2021-10-05 13:39:29.640385",
    "2021-10-05 13:39:29.640389",
    "2021-10-05 13:39:29.640391"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
